# Refresh the cryptos list (Price/Volume(1h) columns, plus two coin-identity
# swaps) with the latest scraped values.
#
# The "Price" column (D) holds plain text even when it looks numeric
# (e.g. "71.94", or thousand-dotted "42.062.00"), matching the source data's
# inlineStr cell type. Writing a numeric-looking string straight to
# .Value lets COM auto-coerce it to a Double, which silently normalizes
# away formatting (e.g. "250.30" -> 250.3, "1.00" -> 1) and changes the
# cell's stored type. So for column D we force text interpretation via
# NumberFormat '@' before the write, then restore the default ("Normal")
# style afterwards so no stray per-cell formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.990.90'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.232.06'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.37%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.633'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '71.94'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.47%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.593'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.10'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +16.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0976'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '58.31'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.29'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.105'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.564.79'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.05'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.79%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.868'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.230.57'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.929.29'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0976'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.24'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.15'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.91'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.15'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +11.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.33%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.56'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +9.13%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.93'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +9.35%  '
$ws.Range('E29').Value = '  +0.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '171.95'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.93'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.12%  '
$ws.Range('E32').Value = '  +3.42%  '
$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.126'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.44%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.63'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0733'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.72'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '26.17'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +20.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.95'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0302'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +14.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.31'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.99'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.44%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '68.37'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.89'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +19.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.208'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +10.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.90'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.81'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.102'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.63%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.69'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +9.92%  '
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('E50').Value = '  +7.56%  '
$ws.Range('E51').Value = '  +2.02%  '
